# Append the 2025-04-27 price row to every sheet in the Solar_Prices workbook.
# Each sheet has a Date/Price table ending at row 56; we add row 57 with the
# new date and a price that repeats the last known price (row 56's value).

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-27"

# Sheet name -> price value to carry forward into the new row.
$sheetPrices = [ordered]@{
    "N-Dense"                   = "38"
    "N-Type"                    = "37.78"
    "N-type Wafer"              = "1.15"
    "Cell Topcon 183mm"         = "0.278"
    "Module Topcon 183mm"       = "0.09"
    "Silver Rear_side"          = "5,424"
    "Silver Busbar front-side"  = "8,121"
    "Silver finger front-side"  = "8,171"
    "USD_CNY"                   = "7.3083"
}

foreach ($sheetName in $sheetPrices.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $price = $sheetPrices[$sheetName]

    $dateCell = $ws.Range("A57")
    $priceCell = $ws.Range("B57")

    # Pre-format the target cells as Text so the numeric/date-looking
    # strings are stored verbatim (matching the rest of the column, which
    # is inline text) instead of being auto-converted to a date serial /
    # number by Excel's normal cell-entry parsing.
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"

    $dateCell.Value = $newDate
    $priceCell.Value = $price

    # Reset to the default style so the new cells don't carry an explicit
    # number format, matching the plain (unstyled) cells elsewhere in the
    # sheet.
    $dateCell.Style = "Normal"
    $priceCell.Style = "Normal"
}
